$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91, shifting existing rows 91-132 down to 92-133.
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new data record.
$ws.Cells.Item(91, 1).Value2 = 9
$ws.Cells.Item(91, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(91, 3).Value2 = "Metropolitana"
$ws.Cells.Item(91, 4).Value2 = 44460
$ws.Cells.Item(91, 5).Value2 = 13
$ws.Cells.Item(91, 6).Value2 = 100112026
$ws.Cells.Item(91, 7).Value2 = "Haba"
$ws.Cells.Item(91, 8).Value2 = "Sin especificar"
$ws.Cells.Item(91, 9).Value2 = "Primera"
$ws.Cells.Item(91, 10).Value2 = 25
$ws.Cells.Item(91, 11).Value2 = 15000
$ws.Cells.Item(91, 12).Value2 = 16000
$ws.Cells.Item(91, 13).Value2 = 15520
$ws.Cells.Item(91, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(91, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(91, 16).Value2 = 621
$ws.Cells.Item(91, 17).Value2 = 25
$ws.Cells.Item(91, 18).Value2 = "Hortaliza"
